$p = $ppt.ActivePresentation

# --- Update the notes master date field (en-IE format: DD/MM/YYYY) ---
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "05/12/2024"

# --- Update the slide master date field (en-US format: MM/DD/YYYY) ---
$master = $p.Designs.Item(1).SlideMaster
$master.Shapes.Item(4).TextFrame.TextRange.Text = "12/5/2024"

# --- Update the date field on every slide layout (en-US format) ---
$layoutDateShapeIndex = @{
    1 = 4
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 4
    12 = 3
    13 = 4
    14 = 4
    15 = 3
    16 = 3
}
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $shIdx = $layoutDateShapeIndex[$li]
    $layout.Shapes.Item($shIdx).TextFrame.TextRange.Text = "12/5/2024"
}

# --- Slide 1 edits ---
$s = $p.Slides.Item(1)

# Title textbox: split "By Qadeer Hussain  Date: 25/10/2024" run into
# "By Qadeer Hussain  Date" and ": 06/12/2024"
$titleShape = $s.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$para2 = $titleTr.Paragraphs(2)
$dateRun = $para2.Runs(2)
$dateRun.Text = "By Qadeer Hussain  Date"
$titleTr.InsertAfter(": 06/12/2024") | Out-Null

# "What are the Core features?" box: reword bullet "Implementing Security"
$featuresShape = $s.Shapes.Item(3)
$featuresTr = $featuresShape.TextFrame.TextRange
$featuresTr.Paragraphs(4).Runs(1).Text = "Security Implementation"
